$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("E29").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
Write-Output "done"
